# Clean-up of input tables (Parameter_HeatingTechnology_Lifetime)
#
# The "min"/"max" lifetime values for the heating-technology table are
# updated:
#   - row 2 (id 11):            min 80 -> 120, max 100 -> 150
#   - rows 3-30 (all other ids): min 30 -> 20,  max 40 -> 30
#
# Headers, ids and the "year" unit column are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First data row (id_heating_technology = 11) gets its own min/max values.
$ws.Range("C2").Value = 120
$ws.Range("D2").Value = 150

# Remaining data rows (3-30) all share the same new min/max values.
$ws.Range("C3:C30").Value = 20
$ws.Range("D3:D30").Value = 30
